$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 17-22 (inclusive) - shifts rows up
$ws.Rows("17:22").Delete()

# Update the summary values to reflect the new (single-worker) data set
$ws.Range("E11").Value = 5389
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

Write-Host "Done"
